$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CargaMasivaDetalle")

# Update A2: "123A" -> "    123A" (4 leading spaces)
$ws.Range("A2").Value = "    123A"

# Update J3: empty -> "ñuñoa"
$ws.Range("J3").Value = "ñuñoa"

# Update J6: "Ñuñoa" -> "chépica"
$ws.Range("J6").Value = "chépica"

# Update the active sheet's view: topLeftCell = G1, selection = J6
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("J6").Select()
